# Handle qualifiers in WQ results field
#
# Adds a new "location" lookup sheet (location_id -> location_name) in
# front of the existing lookup sheets, mirroring the location_id values
# already used on the "eia_area" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Normalize the selection on connectivity_score before we start
#        inserting sheets (best-effort; the sheet otherwise keeps a
#        stale B1 selection left over from a prior edit).
$connectivity = $wb.Worksheets.Item("connectivity_score")
$connectivity.Activate()
$connectivity.Range("A1").Select()

# --- 2. Insert a brand new sheet in front of the first existing sheet
#        and name it "location".
$firstSheet = $wb.Worksheets.Item(1)
$location = $wb.Worksheets.Add($firstSheet)
$location.Name = "location"

# --- 3. Fill in the location_id / location_name lookup table. The
#        "Watts Branch - Upper" value is written before the
#        "location_name" header so new shared-string entries land in
#        the same order the original authoring session produced them.
$location.Range("B28").Value = "Watts Branch - Upper"

$rows = @(
    @("location_id", "location_name"),
    @("RCR01", "Rock Creek Upper"),
    @("RCR09", "Rock Creek Lower"),
    @("TBK01", "Battery Kemble Creek"),
    @("TBR01", "Broad Branch"),
    @("TDA01", "Dalecarlia Tributary"),
    @("TDO01", "Dumbarton Oaks"),
    @("TDU01", "Fort Dupont Tributary"),
    @("TFB01", "Foundry Branch"),
    @("TFC01", "Fort Chaplin Tributary"),
    @("TFD01", "Fort Davis Tributary"),
    @("TFE01", "Fenwick Branch"),
    @("TFS01", "Fort Stanton Tributary"),
    @("THR01", "Hickey Run"),
    @("TKV01", "Klingle Valley Run"),
    @("TLU01", "Luzon Branch"),
    @("TMH01", "Melvin Hazen Valley Branch"),
    @("TNA01", "Nash Run"),
    @("TNS01", "Normanstone"),
    @("TOR01", "Oxon Run"),
    @("TPB01", "Pope Branch"),
    @("TPI01", "Pinehurst Branch"),
    @("TPO01", "Portal Branch"),
    @("TPY01", "Piney Branch"),
    @("TSO01", "Soapstone Creek"),
    @("TTX27", "Texas Avenue Tributary"),
    @("TWB01", "Watts Branch - Lower"),
    @("TWB05", "Watts Branch - Upper")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $location.Cells.Item($r, 1).Value = $rows[$i][0]
    $location.Cells.Item($r, 2).Value = $rows[$i][1]
}

# --- 4. Column widths to match the bestFit sizing used on the other
#        location_id/location_name columns in this workbook.
$location.Columns.Item(1).ColumnWidth = 9.5
$location.Columns.Item(2).ColumnWidth = 24

# --- 5. Selection left on the sheet after data entry.
$location.Range("B32").Select()

# --- 6. The "eia_area" sheet's location_id column widens to match the
#        new "location" sheet's id column (same bestFit target).
$eiaArea = $wb.Worksheets.Item("eia_area")
$eiaArea.Columns.Item(1).ColumnWidth = 9.5

# --- 7. Leave "location" as the active/front-most sheet (tabSelected).
$location.Activate()
